$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet (Property / Value pairs) ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 2.1.0 -> 2.2.0-ballot
$meta.Range("B3").Value = "2.2.0-ballot"

# Date: bump the publication timestamp
$meta.Range("B8").Value = "2025-12-19T09:47:21+00:00"

# Base Definition: pin the FHIR version of the base StructureDefinition
$meta.Range("B18").Value = "http://hl7.org/fhir/StructureDefinition/Extension|4.0.1"

# --- "Elements" sheet (StructureDefinition element table) ---
$elements = $wb.Worksheets.Item("Elements")

# Extension.value[x] Type(s): pin the referenced profile's version
$elements.Range("K6").Value = "Reference(https://interop.esante.gouv.fr/ig/fhir/tddui/StructureDefinition/tddui-questionnaire-response|2.2.0-ballot)`n"

# Widen the "Type(s)" column to fit the now-longer reference text
$elements.Columns.Item(11).ColumnWidth = 90.8
